$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.170.11'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '3.384.32'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.13'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.10'
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.632'
$ws.Range('E7').Value = '  +1.86%  '
$ws.Range('D8').Value = '3.371.90'
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.172'
$ws.Range('E10').Value = '  +5.72%  '
$ws.Range('E11').Value = '  +1.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.74'
$ws.Range('E12').Value = '  -1.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000280'
$ws.Range('E13').Value = '  +2.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.20'
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').Value = '3.910.82'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.31'
$ws.Range('E16').Value = '  -0.62%  '
$ws.Range('D17').Value = '3.383.59'
$ws.Range('E17').Value = '  +0.89%  '
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.91'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = '64.990.35'
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  +1.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '454.43'
$ws.Range('E22').Value = '  -1.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.91'
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.21'
$ws.Range('E24').Value = '  +5.64%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.08'
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('E26').Value = '  +2.27%  '
$ws.Range('E27').Value = '  +0.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.73'
$ws.Range('E28').Value = '  -1.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.73'
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.29'
$ws.Range('E30').Value = '  +3.94%  '
$ws.Range('E31').Value = '  -0.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '63.14'
$ws.Range('E32').Value = '  +7.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.48'
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '577.87'
$ws.Range('E34').Value = '  -0.95%  '
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  +5.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.142'
$ws.Range('E38').Value = '  +0.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.77'
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('D41').Value = '0.0₃0742'
$ws.Range('E41').Value = '  -2.05%  '
$ws.Range('D42').Value = '3.099.44'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0418'
$ws.Range('E43').Value = '  +1.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.78'
$ws.Range('E44').Value = '  -1.21%  '
$ws.Range('E45').Value = '  -1.34%  '
$ws.Range('E46').Value = '  +2.22%  '
$ws.Range('E47').Value = '  -1.25%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '142.09'
$ws.Range('E49').Value = '  +4.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.53'
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.30'
$ws.Range('E51').Value = '  -0.39%  '
